$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff describes a cyclic rotation of the per-row data (columns A, B, E, F, G, H, Q, R)
# across rows 2, 3 and 4:
#   new row2 = old row4
#   new row3 = old row2
#   new row4 = old row3
# Other columns (C, D, I-N, P, S-AY) are identical across the three rows already,
# so only these columns need to change.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values before making any changes.
# Note: use Value2 (not Value) for reliable reads/writes in this COM runtime.
$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range(${col} + "2").Value2
    $orig3[$col] = $ws.Range(${col} + "3").Value2
    $orig4[$col] = $ws.Range(${col} + "4").Value2
}

foreach ($col in $cols) {
    $ws.Range(${col} + "2").Value2 = $orig4[$col]
    $ws.Range(${col} + "3").Value2 = $orig2[$col]
    $ws.Range(${col} + "4").Value2 = $orig3[$col]
}
